# Auto-generated edit script: applies scheduled-runner market-data updates
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 1450.5
$ws.Range("I33").Value = 901
$ws.Range("J33").Value = 2000
$ws.Range("K33").Value = 901
$ws.Range("L33").Value = 2000
$ws.Range("M33").Value = -672
$ws.Range("N33").Value = -2458

$ws.Range("H112").Value = 2204.2222
$ws.Range("I112").Value = 1799.6666
$ws.Range("J112").Value = 2406.5
$ws.Range("K112").Value = 5398.9998
$ws.Range("L112").Value = 7219.5
$ws.Range("M112").Value = -4290.9998
$ws.Range("N112").Value = -9435.5

$ws.Range("H141").Value = 4176
$ws.Range("I141").Value = 4176
$ws.Range("K141").Value = 12528
$ws.Range("M141").Value = -7348

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15646.366
$ws.Range("I32").Value = 12835.429
$ws.Range("K32").Value = 12835.429
$ws.Range("M32").Value = -12548.429

$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()

$ws.Range("H74").Value = 1933.3636
$ws.Range("I74").Value = 1882.8096
$ws.Range("K74").Value = 1882.8096
$ws.Range("M74").Value = -1008.8096

$ws.Range("H77").Value = 1933.3636
$ws.Range("I77").Value = 1882.8096
$ws.Range("K77").Value = 9414.048000000001
$ws.Range("M77").Value = -5046.048000000001

$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H75").Value = 24141.285
$ws.Range("I75").Value = 8164.8335
$ws.Range("K75").Value = 8164.8335
$ws.Range("M75").Value = -7228.8335

$ws.Range("H78").Value = 24141.285
$ws.Range("I78").Value = 8164.8335
$ws.Range("K78").Value = 24494.5005
$ws.Range("M78").Value = -19814.5005

$ws.Range("H105").Value = 3217.7
$ws.Range("I105").Value = 1531.1666
$ws.Range("K105").Value = 1531.1666
$ws.Range("M105").Value = 215.8334

$ws.Range("H109").Value = 79990
$ws.Range("J109").Value = 79990
$ws.Range("L109").Value = 79990
$ws.Range("N109").Value = -82764

$ws.Range("H134").Value = 14698.542
$ws.Range("I134").Value = 13966.429
$ws.Range("K134").Value = 41899.287
$ws.Range("M134").Value = -39364.287

$ws.Range("H138").Value = 49997
$ws.Range("J138").Value = 49997
$ws.Range("L138").Value = 49997
$ws.Range("N138").Value = -60277

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4993.9375
$ws.Range("I31").Value = 4438.25
$ws.Range("J31").Value = 5549.625
$ws.Range("K31").Value = 4438.25
$ws.Range("L31").Value = 5549.625
$ws.Range("M31").Value = -4143.25
$ws.Range("N31").Value = -6139.625

$ws.Range("H34").Value = 4993.9375
$ws.Range("I34").Value = 4438.25
$ws.Range("J34").Value = 5549.625
$ws.Range("K34").Value = 4438.25
$ws.Range("L34").Value = 5549.625
$ws.Range("M34").Value = -4236.25
$ws.Range("N34").Value = -5953.625

$ws.Range("H62").Value = 10643.429
$ws.Range("I62").Value = 11401
$ws.Range("K62").Value = 11401
$ws.Range("M62").Value = -10777

$ws.Range("H65").Value = 10643.429
$ws.Range("I65").Value = 11401
$ws.Range("K65").Value = 57005
$ws.Range("M65").Value = -53885

$ws.Range("H68").Value = 39998.125
$ws.Range("J68").Value = 39998.125
$ws.Range("L68").Value = 39998.125
$ws.Range("N68").Value = -41496.125

$ws.Range("H71").Value = 39998.125
$ws.Range("J71").Value = 39998.125
$ws.Range("L71").Value = 119994.375
$ws.Range("N71").Value = -127482.375

$ws.Range("H74").Value = 39998.125
$ws.Range("J74").Value = 39998.125
$ws.Range("L74").Value = 39998.125
$ws.Range("N74").Value = -41746.125

$ws.Range("H77").Value = 39998.125
$ws.Range("J77").Value = 39998.125
$ws.Range("L77").Value = 119994.375
$ws.Range("N77").Value = -128730.375

$ws.Range("H99").Value = 6566.6665
$ws.Range("I99").Value = 6566.6665
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 6566.6665
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -5068.6665
$ws.Range("N99").ClearContents()

$ws.Range("H107").Value = 875.75
$ws.Range("J107").Value = 1006.5
$ws.Range("L107").Value = 1006.5
$ws.Range("N107").Value = -4846.5

$ws.Range("H122").Value = 2998.7727
$ws.Range("I122").Value = 3262.389
$ws.Range("J122").Value = 1812.5
$ws.Range("K122").Value = 9787.167000000001
$ws.Range("L122").Value = 5437.5
$ws.Range("M122").Value = -7337.167000000001
$ws.Range("N122").Value = -10337.5

$ws.Range("H126").Value = 6566.6665
$ws.Range("I126").Value = 6566.6665
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 19699.9995
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -17229.9995
$ws.Range("N126").ClearContents()

$ws.Range("H132").Value = 1876
$ws.Range("I132").Value = 1872.5
$ws.Range("K132").Value = 5617.5
$ws.Range("M132").Value = -3087.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 16010
$ws.Range("J137").Value = 4000
$ws.Range("L137").Value = 12000
$ws.Range("N137").Value = -22200

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H38").Value = 12500
$ws.Range("I38").Value = 5000
$ws.Range("K38").Value = 5000
$ws.Range("M38").Value = -4537

$ws.Range("H80").Value = 4063
$ws.Range("I80").Value = 3898
$ws.Range("K80").Value = 3898
$ws.Range("M80").Value = -2900

$ws.Range("H83").Value = 4063
$ws.Range("I83").Value = 3898
$ws.Range("K83").Value = 19490
$ws.Range("M83").Value = -14498

$ws.Range("H96").Value = 30237
$ws.Range("I96").Value = 30237
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 30237
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = -27491
$ws.Range("N96").ClearContents()

$ws.Range("H97").Value = 744.2222
$ws.Range("I97").Value = 793.13336
$ws.Range("K97").Value = 793.13336
$ws.Range("M97").Value = -297.13336

$ws.Range("H102").Value = 1561.3636
$ws.Range("J102").Value = 737.5
$ws.Range("L102").Value = 737.5
$ws.Range("N102").Value = -3981.5

$ws.Range("H122").Value = 99429.57000000001
$ws.Range("I122").Value = 117199.8
$ws.Range("K122").Value = 351599.4
$ws.Range("M122").Value = -349149.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5533.8
$ws.Range("I40").Value = 5533.8
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 5533.8
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -5397.8
$ws.Range("N40").ClearContents()

$ws.Range("H55").Value = 1410
$ws.Range("I55").Value = 2625
$ws.Range("J55").Value = 600
$ws.Range("K55").Value = 2625
$ws.Range("L55").Value = 600
$ws.Range("M55").Value = -2452
$ws.Range("N55").Value = -946

$ws.Range("H132").Value = 14226.277
$ws.Range("I132").Value = 13862.286
$ws.Range("J132").Value = 15500.25
$ws.Range("K132").Value = 41586.858
$ws.Range("L132").Value = 46500.75
$ws.Range("M132").Value = -39056.858
$ws.Range("N132").Value = -51560.75

$ws.Range("H136").Value = 5378
$ws.Range("I136").Value = 3503.5
$ws.Range("J136").Value = 7252.5
$ws.Range("K136").Value = 10510.5
$ws.Range("L136").Value = 21757.5
$ws.Range("M136").Value = -7960.5
$ws.Range("N136").Value = -26857.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 5982.6665
$ws.Range("I126").Value = 7299.3335
$ws.Range("K126").Value = 21898.0005
$ws.Range("M126").Value = -19428.0005

$ws.Range("H132").Value = 3133.3333
$ws.Range("I132").Value = 3133.3333
$ws.Range("K132").Value = 9399.999899999999
$ws.Range("M132").Value = -6869.999899999999

$ws.Range("H135").Value = 57617.5
$ws.Range("J135").Value = 57617.5
$ws.Range("L135").Value = 57617.5
$ws.Range("N135").Value = -67757.5
